$wb = $excel.ActiveWorkbook

# --- mon sheet ---
$ws = $wb.Worksheets.Item("mon")
$ws.Range("J2").Value = "CSC442"
$ws.Range("C4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("H15").Value = "CSC425"
$ws.Range("I15").Value = "CSC442"
$ws.Range("G20").Value = "MAT111"
$ws.Range("I20").ClearContents()
$ws.Range("D21").Value = "CST111"
$ws.Range("G21").ClearContents()
$ws.Range("I21").Value = "CHM111"

# --- tue sheet ---
$ws = $wb.Worksheets.Item("tue")
$ws.Range("I3").Value = "CSC424"
$ws.Range("K4").Value = "CSC423"
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E16").Value = "CSC111"
$ws.Range("C20").Value = "MAT111"
$ws.Range("D20").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("G20").ClearContents()
$ws.Range("E21").Value = "GST111"
$ws.Range("F21").ClearContents()
$ws.Range("H21").ClearContents()
$ws.Range("I24").Value = "CSC441"
$ws.Range("J24").Value = "CSC441"

# --- wed sheet ---
$ws = $wb.Worksheets.Item("wed")
$ws.Range("G4").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("G13").Value = "CIS421"
$ws.Range("K13").Value = "CSC111"
$ws.Range("E14").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("D15").ClearContents()
$ws.Range("G20").ClearContents()
$ws.Range("H20").ClearContents()
$ws.Range("I20").Value = "TMC111"
$ws.Range("J21").Value = "MAT112"
$ws.Range("K21").Value = "MAT112"
$ws.Range("F24").Value = "CSC425"
$ws.Range("D25").Value = "CSC424"

# --- thur sheet ---
$ws = $wb.Worksheets.Item("thur")
$ws.Range("E7").Value = "CSC425"
$ws.Range("J7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("G11").Value = "CSC424"
$ws.Range("I14").Value = "CIS421"
$ws.Range("H17").ClearContents()
$ws.Range("B21").Value = "EDS421"
$ws.Range("G21").ClearContents()
$ws.Range("H21").ClearContents()
$ws.Range("I21").ClearContents()
$ws.Range("K24").Value = "CSC442"
$ws.Range("J25").Value = "CSC424"

# --- fri sheet ---
$ws = $wb.Worksheets.Item("fri")
$ws.Range("B2").Value = "BIO111"
$ws.Range("C2").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("E15").Value = "CSC111"
$ws.Range("C17").Value = "BIO111"
$ws.Range("F20").Value = "TMC421"
$ws.Range("G20").ClearContents()
$ws.Range("E21").Value = "GST111"
$ws.Range("C24").Value = "CSC424"
